$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 23.949247360229492
$ws.Range("C2").Value = 6.0804595947265625
$ws.Range("D2").Value = 16.947368621826172
$ws.Range("E2").Value = 57.85714340209961
